$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range of the sheet
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Swap the contents of column A and column B for every row (header + data),
# i.e. change the front-back order of the "category name" / "item name" columns.
for ($r = 1; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $valA = $cellA.Value2
    $valB = $cellB.Value2

    $cellA.Value = $valB
    $cellB.Value = $valA
}
